# ICDC UBC02 changes and MTP 22.11 3 tickets update
# Updates the Neo4j/Cypher queries embedded in the "startup" sheet and adds a
# new "StudyFilesTab" row describing a study-level files query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New / revised query text blocks (single-quoted here-strings so that
# backticks, $ signs and quotes inside the Cypher text are kept verbatim).
# ---------------------------------------------------------------------------

$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (f:file)-[*]->(samp:sample)-->(c)
MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp,demo, c, s, p, diag
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Airedale Terrier', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in [ 'Bladder', 'Bladder, Urethra', 'Bladder, Urethra, Prostate'] and diag.best_response in ['Not Determined']
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$newCasesTabQuery = @'
 MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis) 
 MATCH (samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Airedale Terrier', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in [ 'Bladder', 'Bladder, Urethra', 'Bladder, Urethra, Prostate'] and diag.best_response in ['Not Determined']
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
Order by c.case_id LIMIT 100        
'@

$newFilesTabQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Airedale Terrier', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in [ 'Bladder', 'Bladder, Urethra', 'Bladder, Urethra, Prostate'] and diag.best_response in ['Not Determined']
WITH DISTINCT f, parent, c, demo, diag, s
OPTIONAL MATCH (f)-[*]->(samp:sample)
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN coalesce(f.file_name, '') AS `File Name`, 
 coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`, 
      CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
   coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis 
        Order By f.file_name LIMIT 100
'@

$newStudyFilesQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (sf:file)-->(s)
MATCH (s)<--(c)
MATCH (samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Airedale Terrier', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in [ 'Bladder', 'Bladder, Urethra', 'Bladder, Urethra, Prostate'] and diag.best_response in ['Not Determined']
WITH DISTINCT f,  s, c
WITH
        f, c,  s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c,  s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c,   s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# Trim the trailing newline the here-string literal adds after the last line
# (PowerShell here-strings always end with `\r\n` before the closing `'@`).
$newStatQuery       = $newStatQuery.TrimEnd("`r", "`n")
$newCasesTabQuery   = $newCasesTabQuery.TrimEnd("`r", "`n")
$newFilesTabQuery   = $newFilesTabQuery.TrimEnd("`r", "`n")
$newStudyFilesQuery = $newStudyFilesQuery.TrimEnd("`r", "`n")

# ---------------------------------------------------------------------------
# Write the brand-new shared strings in the same order the original authors'
# edits introduced them, so the shared-string table indexes line up:
#   11 = new StatQuery, 12 = StudyFilesTab, 13 = new CasesTab query,
#   14 = new FilesTab query, 15 = new StudyFiles query
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = $newStatQuery
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B2").Value = $newCasesTabQuery
$ws.Range("B4").Value = $newFilesTabQuery
$ws.Range("B5").Value = $newStudyFilesQuery

# ---------------------------------------------------------------------------
# Row 3 - SamplesTab: query unchanged, but now paired with the new StatQuery
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = $newStatQuery

# ---------------------------------------------------------------------------
# Row 4 - FilesTab: paired with the new StatQuery
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = $newStatQuery

# ---------------------------------------------------------------------------
# Row 5 (new) - StudyFilesTab
# ---------------------------------------------------------------------------
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = $newStatQuery
$ws.Range("C5").WrapText = $true
$ws.Range("D5").Value = $ws.Range("D4").Text
$ws.Range("E5").Value = $ws.Range("E4").Text

# Move the old "blank spacer" row (row 6 in the original layout) down to make
# room for the new StudyFilesTab row, keeping its style.
$ws.Range("C6").WrapText = $true

# ---------------------------------------------------------------------------
# Row heights (auto-fit growth from the longer wrapped Cypher text)
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 390
$ws.Rows.Item(3).RowHeight = 300
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(5).RowHeight = 409.5

# ---------------------------------------------------------------------------
# View state - user had scrolled down to / selected the new row
# ---------------------------------------------------------------------------
$ws.Range("C5").Select()

Write-Output "StudyFilesTab row added; queries updated."
